$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "[0.27660990647713257, 0.3642405133822162]"
$ws.Range("M2").Value = 0.000000000000824229573481716216
$ws.Range("N2").Value = 0.000000000000824229573481716216
$ws.Range("T2").Value = "[0.45837316894349156, 0.5079961967565771]"
